$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.797.89"
$ws.Range("E2").Value = "  +1.24%  "

$ws.Range("D3").Value = "2.436.88"
$ws.Range("E3").Value = "  +0.89%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "570.04"
$ws.Range("E5").Value = "  +1.25%  "

$ws.Range("D6").Value = "146.28"
$ws.Range("E6").Value = "  +2.64%  "

$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.17%  "

$ws.Range("D8").Value = "0.534"
$ws.Range("E8").Value = "  +0.64%  "

$ws.Range("D9").Value = "0.112"
$ws.Range("E9").Value = "  +1.76%  "

$ws.Range("E10").Value = "  +0.28%  "

$ws.Range("D11").Value = "5.32"
$ws.Range("E11").Value = "  +1.87%  "

$ws.Range("D12").Value = "0.358"
$ws.Range("E12").Value = "  +2.28%  "

$ws.Range("D13").Value = "26.85"
$ws.Range("E13").Value = "  +4.92%  "

$ws.Range("E14").Value = "  +4.53%  "

$ws.Range("D15").Value = "2.877.04"
$ws.Range("E15").Value = "  +0.77%  "

$ws.Range("D16").Value = "62.645.48"
$ws.Range("E16").Value = "  +1.17%  "

$ws.Range("D17").Value = "2.429.61"
$ws.Range("E17").Value = "  +0.65%  "

$ws.Range("D18").Value = "11.27"
$ws.Range("E18").Value = "  +0.34%  "

$ws.Range("D19").Value = "7.06"
$ws.Range("E19").Value = "  +3.52%  "

$ws.Range("D20").Value = "325.07"
$ws.Range("E20").Value = "  +1.20%  "

$ws.Range("D21").Value = "4.19"
$ws.Range("E21").Value = "  +1.37%  "

$ws.Range("D22").Value = "1.05"
$ws.Range("E22").Value = "  +4.90%  "

$ws.Range("E23").Value = "  +4.91%  "

$ws.Range("D24").Value = "67.24"

$ws.Range("D25").Value = "618.64"
$ws.Range("E25").Value = "  +9.50%  "

$ws.Range("D26").Value = "8.68"
$ws.Range("E26").Value = "  +0.47%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000102"
$ws.Range("E27").Value = "  +9.17%  "

$ws.Range("D28").Value = "2.556.96"
$ws.Range("E28").Value = "  +0.91%  "

$ws.Range("D29").Value = "8.48"
$ws.Range("E29").Value = "  +3.67%  "

$ws.Range("E30").Value = "  -0.03%  "

$ws.Range("D31").Value = "1.47"
$ws.Range("E31").Value = "  +4.69%  "

$ws.Range("E32").Value = "  -3.23%  "

$ws.Range("E33").Value = "  +0.31%  "

$ws.Range("E34").Value = "  -0.30%  "

$ws.Range("E35").Value = "  +3.18%  "

$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.05%  "

$ws.Range("E37").Value = "  +0.96%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").Value = "5.41"
$ws.Range("E38").Value = "  -0.52%  "

$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").Value = "18.77"
$ws.Range("E39").Value = "  +1.16%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "1.83"
$ws.Range("E40").Value = "  +1.80%  "

$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "148.25"
$ws.Range("E41").Value = "  -3.18%  "

$ws.Range("E42").Value = "  +13.92%  "

$ws.Range("E43").Value = "  +0.26%  "

$ws.Range("D44").Value = "149.83"
$ws.Range("E44").Value = "  +0.66%  "

$ws.Range("E45").Value = "  +2.31%  "

$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0540"
$ws.Range("E46").Value = "  +1.72%  "

$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "20.74"
$ws.Range("E47").Value = "  +4.50%  "

$ws.Range("D48").Value = "0.605"
$ws.Range("E48").Value = "  +1.93%  "

$ws.Range("D49").Value = "0.0233"
$ws.Range("E49").Value = "  +2.91%  "

$ws.Range("D50").Value = "0.0923"
$ws.Range("E50").Value = "  +0.28%  "

$ws.Range("E51").Value = "  +4.49%  "
